# Weekly roll of the "Perejil" price series:
# - a new most-recent fortnight (date 44523) is inserted at the top of the
#   date-block (rows 84:85)
# - every existing fortnight block shifts down by one pair of rows (84:85 -> 86:87, etc.)
# - the oldest block (previously rows 114:115, date 44400) falls off the
#   bottom of the shifted block and is appended as new rows 116:117

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the last existing pair of rows (114:115) onto new rows 116:117
#    before anything else touches them.
$ws.Range("A114:R115").Copy($ws.Range("A116:R117"))

# 2) Shift the "Fecha" (date) column down by one pair of rows for the block
#    that used to run rows 84:113 (now landing on 86:115). Copy from the
#    bottom up isn't required since Copy reads the source before writing,
#    but to be safe or order-independent we copy in one shot.
$ws.Range("D84:D113").Copy($ws.Range("D86:D115"))

# 3) Write the brand-new date for the newly inserted top pair of rows.
$ws.Range("D84").Value = 44523
$ws.Range("D85").Value = 44523
